$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Major Updates - Scenario, etc." --------------------------------
# The mev_type lookup table gains a "description" and "country" column,
# the M1 row is replaced by a PSR (Personal Savings Rate) row, a new
# UNRATE (Unemployment Rate) row is inserted, and PDI's type flips from
# "rate" to "level".

# Insert a new row at position 5; this shifts the existing CPI row from
# row 5 down to row 6, making room for the new UNRATE row.
$ws.Rows("5").Insert()

# PDI's "type" changes from "rate" to "level".
$ws.Range("B4").Value = "level"

# New column C: human-readable description for each series.
$ws.Range("C1").Value = "description"
$ws.Range("C2").Value = "Nominal GDP"
$ws.Range("C3").Value = "Personal Savings Rate"
$ws.Range("C4").Value = "Personal Disposable Income"
$ws.Range("C5").Value = "Unemployment Rate"
$ws.Range("C6").Value = "Consumer Price Index"

# New column D: country code for each series.
$ws.Range("D1").Value = "country"
$ws.Range("D2").Value = "US"
$ws.Range("D3").Value = "US"
$ws.Range("D4").Value = "US"
$ws.Range("D5").Value = "US"
$ws.Range("D6").Value = "US"

# Row 3: the old "M1" code/row is replaced by "PSR", whose type is "rate".
$ws.Range("A3").Value = "PSR"
$ws.Range("B3").Value = "rate"

# Row 5 (newly inserted): the Unemployment Rate series.
$ws.Range("A5").Value = "UNRATE"
$ws.Range("B5").Value = "rate"

# Auto-fit the new description column to its contents.
$ws.Columns("C").AutoFit()

# Page setup / selection bookkeeping.
$ws.PageSetup.Orientation = 1
[void]$ws.Range("H5").Select()
